# Applies the "Added - Missing Examples in Aspose.Slides" edit:
#  1. Update the cached text of the date placeholder on the (only) slide
#     layout from 11/7/2009 to 5/26/2016.
#  2. Add two slide guides (horizontal @ 2160, vertical @ 2880) to the
#     presentation - best effort (host may not persist these).
#  3. Remove the Aspose "Evaluation only." watermark TextBox from slide 1.
#  4. On slide 2, shrink/move the old watermark TextBox and blank out its
#     text, then add a small new empty "TextBox 1" shape next to it
#     (mirrors what Aspose.Slides' own sample now looks like once the
#     evaluation watermark was removed from the library).

$p = $ppt.ActivePresentation

# --- 1. Date placeholder text cached on the slide layout -------------------
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(1)
$dateShape = $layout.Shapes.Item("Date Placeholder 1")
$dateShape.TextFrame.TextRange.Text = "5/26/2016"

# --- 2. Slide guides (best effort; harmless if host ignores it) ------------
$horizGuide = $p.Guides.Add(1, 2160)
$vertGuide = $p.Guides.Add(2, 2880)

# --- 3. Slide 1: delete the Aspose evaluation watermark textbox ------------
$slide1 = $p.Slides.Item(1)
$watermark1 = $slide1.Shapes.Item("TextBox")
$watermark1.Delete()

# --- 4. Slide 2: repurpose the watermark textbox + add a new textbox -------
$slide2 = $p.Slides.Item(2)
$watermark2 = $slide2.Shapes.Item("TextBox")
$watermark2.Left = 352.743408203125
$watermark2.Top = 240.2852783203125
$watermark2.Width = 14.545669555664062
$watermark2.Height = 59.42472457885742
$watermark2.TextFrame.TextRange.Text = ""

$newBox = $slide2.Shapes.AddTextbox(1, 291.9609680175781, 230.31056213378906, 14.545748710632324, 29.081260681152344)
$newBox.Name = "TextBox 1"
$newBox.Fill.Visible = $false
$newBox.TextFrame.WordWrap = $false
$newBox.TextFrame.AutoSize = 1
$newBox.TextFrame.TextRange.Text = ""
